$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new block (rows 39-47) is a copy of the existing "Rules String Hello2"
# block (rows 27-35), renamed to "Rules String Hello3 (Integer hour)".
# Copy the whole formatted range (values + styles + merged header cell) in
# one shot, then just retarget the header text.
$src = $ws.Range("B27:E35")
$dst = $ws.Range("B39")
$src.Copy($dst)

$ws.Range("B39").Value = "Rules String Hello3 (Integer hour)"

# Comments that document the new Decision Table, mirroring the ones already
# present on the other tables (B3/B4/C4/E4/C5/E5, B16, B27/B28/C28/E28/C29/E29).
$ws.Range("B39").AddComment("This is so-called Decision Table Header. It starts with the keyword ""Rules"".")

$ws.Range("B40").AddComment("`nRule column header. Rule column is used to to name particular rule rows for documentation and tracing purposes. It is also useful to create rule rows that span more than one cell vertically (this will be explained in one of the next tutorials)`n")

$ws.Range("C40").AddComment("Condition column header. Must start with ""C""")

$ws.Range("E40").AddComment("Return column header. Must start with ""RET"".  ")

$ws.Range("C41").AddComment("Condition expression. Must have type boolean. As you can see condition uses parameter hour from Method Header and variable min that defines column data. When condition is evaluated for each row, the cell value from this row is assigned to variable min")

$ws.Range("E41").AddComment("This is return expression performed for the first row where all conditions have been satisfied. The variable greeting is substittuted with a cell value from the rule row")
